$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.008.91"
$ws.Range("E2").Value = "  +5.72%  "
$ws.Range("D3").Value = "3.713.48"
$ws.Range("E3").Value = "  +19.55%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'617.93"
$ws.Range("E5").Value = "  +7.69%  "
$ws.Range("D6").Value = "'182.69"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("D7").Value = "3.710.59"
$ws.Range("E7").Value = "  +19.51%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.542"
$ws.Range("E9").Value = "  +5.60%  "
$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = "  +7.73%  "
$ws.Range("D11").Value = "'6.60"
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("D12").Value = "'0.502"
$ws.Range("E12").Value = "  +6.95%  "
$ws.Range("D13").Value = "'40.52"
$ws.Range("E13").Value = "  +12.31%  "
$ws.Range("E14").Value = "  +6.17%  "
$ws.Range("D15").Value = "4.332.85"
$ws.Range("E15").Value = "  +19.45%  "
$ws.Range("D16").Value = "3.721.76"
$ws.Range("E16").Value = "  +19.79%  "
$ws.Range("D17").Value = "71.060.48"
$ws.Range("E17").Value = "  +5.92%  "
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").Value = "'7.53"
$ws.Range("E19").Value = "  +6.80%  "
$ws.Range("D20").Value = "'518.82"
$ws.Range("E20").Value = "  +5.30%  "
$ws.Range("D21").Value = "'16.88"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").Value = "'9.36"
$ws.Range("E22").Value = "  +20.30%  "
$ws.Range("D23").Value = "'0.745"
$ws.Range("E23").Value = "  +8.47%  "
$ws.Range("E24").Value = "  +12.81%  "
$ws.Range("D25").Value = "'88.74"
$ws.Range("E25").Value = "  +6.13%  "
$ws.Range("D26").Value = "'13.55"
$ws.Range("E26").Value = "  +7.54%  "
$ws.Range("D27").Value = "'11.22"
$ws.Range("E27").Value = "  +11.29%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").Value = "'2.54"
$ws.Range("E29").Value = "  +9.80%  "
$ws.Range("D30").Value = "'8.19"
$ws.Range("E30").Value = "  +3.27%  "
$ws.Range("D31").Value = "'2.89"
$ws.Range("E31").Value = "  +11.18%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "'0.0000112"
$ws.Range("E32").Value = "  +18.94%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'31.88"
$ws.Range("E33").Value = "  +13.40%  "
$ws.Range("E34").Value = "  +4.36%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").Value = "'1.04"
$ws.Range("E36").Value = "  +10.57%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "'6.13"
$ws.Range("E37").Value = "  +9.96%  "
$ws.Range("D38").Value = "'2.25"
$ws.Range("E38").Value = "  +12.01%  "
$ws.Range("D39").Value = "'0.346"
$ws.Range("E39").Value = "  +10.84%  "
$ws.Range("E40").Value = "  +8.79%  "
$ws.Range("D41").Value = "'51.62"
$ws.Range("E41").Value = "  +5.17%  "
$ws.Range("D42").Value = "'436.16"
$ws.Range("E42").Value = "  +17.65%  "
$ws.Range("D43").Value = "'44.93"
$ws.Range("E43").Value = "  -5.73%  "
$ws.Range("D44").Value = "3.151.70"
$ws.Range("E44").Value = "  +12.78%  "
$ws.Range("D45").Value = "'8.84"
$ws.Range("E45").Value = "  +6.54%  "
$ws.Range("D46").Value = "'2.87"
$ws.Range("E46").Value = "  +5.18%  "
$ws.Range("D47").Value = "'0.0368"
$ws.Range("E47").Value = "  +6.71%  "
$ws.Range("D48").Value = "'28.20"
$ws.Range("E48").Value = "  +10.84%  "
$ws.Range("D49").Value = "'140.73"
$ws.Range("E49").Value = "  +3.59%  "
$ws.Range("D51").Value = "'2.48"
$ws.Range("E51").Value = "  +8.82%  "
